$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.485.87'
$ws.Range('E2').Value = '  +2.72%  '
$ws.Range('D3').Value = '2.661.78'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.07'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.63'
$ws.Range('E6').Value = '  +5.07%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('E9').Value = '  +8.67%  '
$ws.Range('E10').Value = '  +5.24%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.68'
$ws.Range('E13').Value = '  +6.02%  '
$ws.Range('E14').Value = '  +16.69%  '
$ws.Range('D15').Value = '3.141.77'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '65.259.16'
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').Value = '2.662.62'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.80'
$ws.Range('E18').Value = '  +4.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.94'
$ws.Range('E19').Value = '  +3.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '361.54'
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.40'
$ws.Range('E21').Value = '  +6.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.51'
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.73'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.65'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000104'
$ws.Range('E26').Value = '  +18.35%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.67'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.35'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.166'
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '556.29'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').Value = '  +8.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.71'
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  +5.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.434'
$ws.Range('E36').Value = '  +3.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.69'
$ws.Range('E37').Value = '  +5.31%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.43'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.02'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.63'
$ws.Range('E42').Value = '  +7.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '168.36'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.22'
$ws.Range('E44').Value = '  +3.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0625'
$ws.Range('E45').Value = '  +6.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.34'
$ws.Range('E46').Value = '  +8.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.31'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E48').Value = '  +4.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0266'
$ws.Range('E49').Value = '  +5.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0989'
$ws.Range('E50').Value = '  +2.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.85'
$ws.Range('E51').Value = '  +2.34%  '
